# Append a fresh scrape (2025-10-09 06:27 JST) to the top of the "ランサーズ"
# list and trim the sheet back down to just the two surviving rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Drop the stale tail rows (old rows 4-18) -----------------------------
$ws.Rows("4:18").Delete()

# --- Row 2: "海外仕入れ元サイト..." (was row 5 in the previous scrape) ----
$ws.Range("A2").Value = "2025-10-09 06:27:33"
$ws.Range("B2").Value = "海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)"
$ws.Range("D2").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5251319"
$ws.Range("G2").Value = 135
$ws.Range("H2").Value = "◆ツール,スクレイピング ◇サイト"

# --- Row 3: new entry "【急募】クローン作成アプリ開発" -------------------
$ws.Range("A3").Value = "2025-10-09 06:27:33"
$ws.Range("B3").Value = "【急募】クローン作成アプリ開発"
$ws.Range("D3").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5409967"
$ws.Range("G3").Value = 93
$ws.Range("H3").Value = "◆開発 ◇アプリ"

# --- Hyperlinks: only F2 / F3 should remain, pointing at the new URLs ----
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5251319")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5409967")
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"

# --- Column widths: B 52 -> 51, D 41 -> 28 --------------------------------
$ws.Columns("B").ColumnWidth = 50.166666
$ws.Columns("D").ColumnWidth = 27.166666
